$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("工作表1")

# Q4: set date value
$ws.Range("Q4").Value = 45861

# Row 7: add A7 value, change I7 value
$ws.Range("A7").Value = 45867
$ws.Range("I7").Value = 45826

# Row 8: change I8 value
$ws.Range("I8").Value = 45847

# Row 9: add I9 value
$ws.Range("I9").Value = 45861

# Window / view state changes
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = $ws.Range("G1").Column
$ws.Range("U8").Select()
